# Insert a new weekly price record at row 408 for
# "Macroferia Regional de Talca" / Piña / Caramelo / Segunda.
# This pushes the existing rows 408:500 down to 409:501 (dimension grows
# from A1:T500 to A1:T501) and populates the newly opened row 408 with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(408).Insert()

$ws.Range("A408").Value = 5
$ws.Range("B408").Value = "Macroferia Regional de Talca"
$ws.Range("C408").Value = "Maule"
$ws.Range("D408").Value = 45275
$ws.Range("E408").Value = 7
$ws.Range("F408").Value = "Fruta"
$ws.Range("G408").Value = 100108
$ws.Range("H408").Value = "Tropicales y subtropicales"
$ws.Range("I408").Value = 100108005
$ws.Range("J408").Value = "Piña"
$ws.Range("K408").Value = "Caramelo"
$ws.Range("L408").Value = "Segunda"
$ws.Range("M408").Value = 250
$ws.Range("N408").Value = 22000
$ws.Range("O408").Value = 22000
$ws.Range("P408").Value = 22000
$ws.Range("Q408").Value = "$/caja 14 unidades"
$ws.Range("R408").Value = "Ecuador"
$ws.Range("S408").Value = 1571
$ws.Range("T408").Value = 14
